# Insert a new data row before the current row 235 ("Camote" / 44566 / ...),
# which pushes that row and everything below it down by one (to 236..336).
# Then populate the newly inserted row 235 with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(235).Insert()

$ws.Cells.Item(235,1).Value  = 5
$ws.Cells.Item(235,2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(235,3).Value  = "Maule"
$ws.Cells.Item(235,4).Value  = 44839
$ws.Cells.Item(235,5).Value  = 7
$ws.Cells.Item(235,6).Value  = 100112045
$ws.Cells.Item(235,7).Value  = "Zapallo"
$ws.Cells.Item(235,8).Value  = "Paine"
$ws.Cells.Item(235,9).Value  = "1a (guarda)"
$ws.Cells.Item(235,10).Value = 1600
$ws.Cells.Item(235,11).Value = 400
$ws.Cells.Item(235,12).Value = 400
$ws.Cells.Item(235,13).Value = 400
$ws.Cells.Item(235,14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(235,15).Value = "Región del Maule"
$ws.Cells.Item(235,16).Value = 400
$ws.Cells.Item(235,17).Value = 1
$ws.Cells.Item(235,18).Value = "Hortaliza"
